$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column A (factor labels)
$colA = @("f1", "f2", "f3", "f4", "f1", "f2", "f3", "f4", "f1", "f2", "f3", "f4")

# New values for column B (soil type labels)
$colB = @("arcilloso", "arcilloso", "arcilloso", "arcilloso", "arenoso", "arenoso", "arenoso", "arenoso", "franco arenoso", "franco arenoso", "franco arenoso", "franco arenoso")

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
}

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

# Autofit column B to match the new text content
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Update the active selection
$ws.Range("E4").Select() | Out-Null
